$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows appended after the existing last row (226).
# Columns: A = date (serial), B = nuovi pos., C = somma mobile 7gg., D = somma mobile 7gg. per 100mila abitanti
$data = @(
    @(44301, 7, 38, 350.0690925840626),
    @(44302, 9, 46, 423.7678489175495),
    @(44303, 6, 41, 377.7061262091202)
)

$startRow = 227
$formatSource = $ws.Cells.Item($startRow - 1, 1)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $row[0]
    $formatSource.Copy()
    $cellA.PasteSpecial(-4122) # xlPasteFormats, matches style of the preceding date cell (s="2")

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
